# Issue#17  REQ 3.1  ADD GREEN UNIT TEST for scenario in which SAP generates additional column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column L, using the header style (s=3) seen on the header row
$ws.Cells.Item(1, 12).Value = "dupa"

# Fill column L for data rows 2-26 with the same "dupa" value/style
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 12).Value = "dupa"
    $ws.Cells.Item($r, 12).WrapText = $true
}

# Move the active selection to M20, matching the post-edit sheet view
$ws.Range("M20").Select()
